$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows: one between the current rows 2 and 3 (new row 3),
# and one between the current rows 4 and 5, i.e. after the row-3 insert
# shifts everything down, at row 5. This turns the original 4 data rows
# (timesteps 0,1,2) into 5 rows with interpolated timesteps 1 and 3 added
# in between, matching "extra timesteps" in the commit message.
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(5).Insert()

# Copy the formatting (font/border/alignment) used for column A down into
# the two freshly inserted rows so they keep the same bold/centered/bordered
# look as the rest of column A.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A5").PasteSpecial(-4122)

# Re-number the timestep index column sequentially: 0,1,2,3,4
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4

# The newly inserted rows inherited the special "date/number" style that was
# on E/G of the row above (row 2) via Rows.Insert; the target rows should be
# plain (no special style) on those columns, so clear that formatting back off.
$ws.Range("E3").ClearFormats()
$ws.Range("G3").ClearFormats()
$ws.Range("E5").ClearFormats()
$ws.Range("G5").ClearFormats()

# Fill the new row 3 with formulas that average the rows immediately above
# and below it (rows 2 and 4).
$ws.Range("B3").Formula = "=(B2+B4)/2"
$ws.Range("C3:H3").Formula = "=(C2+C4)/2"

# Fill the new row 5 with formulas that average the rows immediately above
# and below it (rows 4 and 6).
$ws.Range("B5").Formula = "=(B4+B6)/2"
$ws.Range("C5:H5").Formula = "=(C4+C6)/2"

# Match the saved selection/active cell shown in the updated workbook.
$null = $ws.Range("F10").Select()
